$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1904761904761905
$ws.Range("C2").Value = 0.5746031746031746
$ws.Range("J2").Value = 0.009523809523809525
$ws.Range("P2").Value = 0.1365079365079365
$ws.Range("S2").Value = 0.08888888888888889
$ws.Range("B3").Value = 0.005076142131979695
$ws.Range("C3").Value = 0.05076142131979695
$ws.Range("J3").Value = 0.04060913705583756
$ws.Range("P3").Value = 0.7258883248730964
$ws.Range("S3").Value = 0.1776649746192893
$ws.Range("J4").Value = 0.04444444444444445
$ws.Range("P4").Value = 0.7777777777777778
$ws.Range("S4").Value = 0.1777777777777778
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.0502283105022831
$ws.Range("E6").Value = 0.0091324200913242
$ws.Range("F6").Value = 0.0547945205479452
$ws.Range("J6").Value = 0.2465753424657534
$ws.Range("O6").Value = 0.0136986301369863
$ws.Range("Q6").Value = 0.2146118721461187
$ws.Range("R6").Value = 0.0684931506849315
$ws.Range("S6").Value = 0.3424657534246575
$ws.Range("B7").Value = 0.1161290322580645
$ws.Range("D7").Value = 0.02580645161290323
$ws.Range("E7").Value = 0.006451612903225806
$ws.Range("F7").Value = 0.06451612903225806
$ws.Range("J7").Value = 0.1096774193548387
$ws.Range("O7").Value = 0.03870967741935484
$ws.Range("Q7").Value = 0.1548387096774194
$ws.Range("R7").Value = 0.08387096774193549
$ws.Range("S7").Value = 0.4
$ws.Range("B8").Value = 0.1243386243386243
$ws.Range("D8").Value = 0.01587301587301587
$ws.Range("F8").Value = 0.05026455026455026
$ws.Range("J8").Value = 0.1137566137566138
$ws.Range("O8").Value = 0.01587301587301587
$ws.Range("Q8").Value = 0.2116402116402116
$ws.Range("R8").Value = 0.08201058201058201
$ws.Range("S8").Value = 0.3862433862433862
$ws.Range("B9").Value = 0.1012658227848101
$ws.Range("D9").Value = 0.02531645569620253
$ws.Range("F9").Value = 0.08860759493670886
$ws.Range("J9").Value = 0.1308016877637131
$ws.Range("O9").Value = 0.02109704641350211
$ws.Range("Q9").Value = 0.2109704641350211
$ws.Range("R9").Value = 0.08016877637130802
$ws.Range("S9").Value = 0.3417721518987342
$ws.Range("B10").Value = 0.1200657894736842
$ws.Range("D10").Value = 0.02384868421052632
$ws.Range("E10").Value = 0.0008223684210526315
$ws.Range("F10").Value = 0.07648026315789473
$ws.Range("J10").Value = 0.1208881578947368
$ws.Range("O10").Value = 0.02138157894736842
$ws.Range("Q10").Value = 0.1990131578947368
$ws.Range("R10").Value = 0.09375
$ws.Range("S10").Value = 0.34375
$ws.Range("G11").Value = 0.1777777777777778
$ws.Range("J11").Value = 0.1185185185185185
$ws.Range("K11").Value = 0.2296296296296296
$ws.Range("L11").Value = 0.4666666666666667
$ws.Range("S11").Value = 0.007407407407407408
$ws.Range("G12").Value = 0.7076923076923077
$ws.Range("J12").Value = 0.2384615384615385
$ws.Range("K12").Value = 0.007692307692307693
$ws.Range("L12").Value = 0.02307692307692308
$ws.Range("S12").Value = 0.02307692307692308
$ws.Range("G13").Value = 0.6896551724137931
$ws.Range("J13").Value = 0.2758620689655172
$ws.Range("S13").Value = 0.03448275862068965
$ws.Range("F15").Value = 0.02788844621513944
$ws.Range("H15").Value = 0.1155378486055777
$ws.Range("I15").Value = 0.08764940239043825
$ws.Range("J15").Value = 0.4143426294820717
$ws.Range("K15").Value = 0.0398406374501992
$ws.Range("M15").Value = 0.0199203187250996
$ws.Range("O15").Value = 0.04780876494023904
$ws.Range("S15").Value = 0.2470119521912351
$ws.Range("F16").Value = 0.02336448598130841
$ws.Range("H16").Value = 0.1401869158878505
$ws.Range("I16").Value = 0.1308411214953271
$ws.Range("J16").Value = 0.3457943925233645
$ws.Range("K16").Value = 0.09813084112149532
$ws.Range("M16").Value = 0.01869158878504673
$ws.Range("O16").Value = 0.07476635514018691
$ws.Range("S16").Value = 0.1682242990654206
$ws.Range("F17").Value = 0.0138568129330254
$ws.Range("H17").Value = 0.1778290993071593
$ws.Range("I17").Value = 0.1478060046189376
$ws.Range("J17").Value = 0.4110854503464203
$ws.Range("K17").Value = 0.07852193995381063
$ws.Range("M17").Value = 0.009237875288683603
$ws.Range("O17").Value = 0.07390300230946882
$ws.Range("S17").Value = 0.08775981524249422
$ws.Range("F18").Value = 0.02590673575129534
$ws.Range("H18").Value = 0.1398963730569948
$ws.Range("I18").Value = 0.1709844559585492
$ws.Range("J18").Value = 0.3678756476683938
$ws.Range("K18").Value = 0.1036269430051813
$ws.Range("M18").Value = 0.005181347150259068
$ws.Range("O18").Value = 0.08290155440414508
$ws.Range("S18").Value = 0.1036269430051813
$ws.Range("F19").Value = 0.01571164510166359
$ws.Range("H19").Value = 0.2042513863216266
$ws.Range("I19").Value = 0.08317929759704251
$ws.Range("J19").Value = 0.3964879852125693
$ws.Range("K19").Value = 0.1081330868761553
$ws.Range("M19").Value = 0.01386321626617375
$ws.Range("O19").Value = 0.09611829944547134
$ws.Range("S19").Value = 0.0822550831792976
